# Update "想去人数" (interest/visit counter) values on the 展览, 演出 and
# 全部类型 sheets to reflect the latest generated data (gh-pages output at
# commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15161
$ws1.Range("F3").Value = 19498
$ws1.Range("F14").Value = 206
$ws1.Range("F22").Value = 8168
$ws1.Range("F28").Value = 13
$ws1.Range("F31").Value = 6524
$ws1.Range("F34").Value = 185
$ws1.Range("F37").Value = 5554

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 25

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15161
$ws4.Range("F3").Value = 19498
$ws4.Range("F14").Value = 206
$ws4.Range("F23").Value = 8168
$ws4.Range("F29").Value = 13
$ws4.Range("F32").Value = 25
$ws4.Range("F34").Value = 6524
$ws4.Range("F37").Value = 185
$ws4.Range("F40").Value = 5554
